# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 23, shifting the existing
# rows 23:85 down to 24:86 (matches the target dimension A1:R86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; Excel shifts rows 23:85 -> 24:86
# and copies the formatting (incl. the date style on column D) from
# the row above, same as a manual "Insert Sheet Rows" in the UI.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record.
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44838
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100114001
$ws.Range("G23").Value = "Papa"
$ws.Range("H23").Value = "Cardinal"
$ws.Range("I23").Value = "1a (cosecha)"
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 15500
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Provincia de Melipilla"
$ws.Range("P23").Value = 620
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
